$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are treated as Text so values like "1.000" or
# "24.623.86" are not auto-converted/normalized into numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.623.86'
$ws.Range("E2").Value = '  +3.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.696.66'
$ws.Range("E3").Value = '  +1.86%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.09'
$ws.Range("E5").Value = '  +2.19%  '

$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("E7").Value = '  +1.73%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4017'
$ws.Range("E8").Value = '  +1.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.525'
$ws.Range("E9").Value = '  +4.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9998'
$ws.Range("E10").Value = '  +0.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.61'
$ws.Range("E11").Value = '  -2.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08776'
$ws.Range("E12").Value = '  +1.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.242'
$ws.Range("E13").Value = '  +6.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.30'
$ws.Range("E14").Value = '  +2.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.191'
$ws.Range("E15").Value = '  +11.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001314'
$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.697.48'
$ws.Range("E17").Value = '  +1.86%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '99.92'
$ws.Range("E18").Value = '  +0.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07077'
$ws.Range("E19").Value = '  +2.72%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.69'
$ws.Range("E20").Value = '  +2.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.070'
$ws.Range("E21").Value = '  +6.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.23'
$ws.Range("E23").Value = '  +2.96%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.623.36'
$ws.Range("E24").Value = '  +3.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.123'
$ws.Range("E25").Value = '  +9.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.337'
$ws.Range("E26").Value = '  +1.22%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.86'
$ws.Range("E27").Value = '  +5.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.90'
$ws.Range("E28").Value = '  +1.23%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '136.61'
$ws.Range("E29").Value = '  +4.78%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.199'
$ws.Range("E30").Value = '  +1.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.507'
$ws.Range("E31").Value = '  +9.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.884.79'
$ws.Range("E32").Value = '  +1.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.080'
$ws.Range("E33").Value = '  -3.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08605'
$ws.Range("E34").Value = '  +0.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.144'
$ws.Range("E35").Value = '  +5.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.55'
$ws.Range("E36").Value = '  +9.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2739'
$ws.Range("E37").Value = '  +3.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.922'
$ws.Range("E38").Value = '  +0.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.44'
$ws.Range("E39").Value = '  -0.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09138'
$ws.Range("E40").Value = '  +3.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02732'
$ws.Range("E41").Value = '  +7.58%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.491'
$ws.Range("E42").Value = '  +2.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7651'
$ws.Range("E43").Value = '  +0.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7173'
$ws.Range("E44").Value = '  +1.55%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.579'
$ws.Range("E45").Value = '  +6.68%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.60'
$ws.Range("E46").Value = '  +3.85%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.223'
$ws.Range("E47").Value = '  +2.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  +0.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.03'
$ws.Range("E49").Value = '  +0.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.319'
$ws.Range("E50").Value = '  +8.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07987'
$ws.Range("E51").Value = '  +2.19%  '
